$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the boolean formulas in E2 and E3 with the literal text "TRUE"
# (the source data used =TRUE() which stored a numeric boolean; the fix
# stores the word "TRUE" as a plain text value instead).
$ws.Range("E2").Formula = "=""TRUE"""
$ws.Range("E2").Copy()
$ws.Range("E2").PasteSpecial(-4163)

$ws.Range("E3").Formula = "=""TRUE"""
$ws.Range("E3").Copy()
$ws.Range("E3").PasteSpecial(-4163)

$excel.CutCopyMode = $false

# Update the selection to match the new active range
$ws.Range("E2:E3").Select()
